# save politik untuk sharing ke sma 1
# Update tabulasi figures on "tabulasi all" sheet with refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Real Number Jenis Kelamin (first block) ---
$ws.Range("C3").Value = 242
$ws.Range("C4").Value = 258

# --- Percentage Jenis Kelamin ---
$ws.Range("C9").Value = 48.4
$ws.Range("C10").Value = 51.6

# --- Real Number and Percentage Jenis Kelamin (gender block) ---
$ws.Range("C15").Value = 242
$ws.Range("C16").Value = 258
$ws.Range("C18").Value = 48.4
$ws.Range("C19").Value = 51.6

# --- Real Number and Percentage Jenis Kelamin (usia block) ---
$ws.Range("C23").Value = 242
$ws.Range("C24").Value = 258
$ws.Range("C26").Value = 18
$ws.Range("C27").Value = 43
$ws.Range("C28").Value = 66
$ws.Range("C29").Value = 39
$ws.Range("C30").Value = 334

# --- Crosstabulasi gender vs usia ---
$ws.Range("C36").Value = 2.89256198347107
$ws.Range("D36").Value = 4.26356589147287
$ws.Range("E36").Value = 3.6

$ws.Range("C37").Value = 9.91735537190083
$ws.Range("D37").Value = 7.36434108527132
$ws.Range("E37").Value = 8.6

$ws.Range("C38").Value = 11.1570247933884
$ws.Range("D38").Value = 15.1162790697674
$ws.Range("E38").Value = 13.2

$ws.Range("C39").Value = 9.91735537190083
$ws.Range("D39").Value = 5.81395348837209
$ws.Range("E39").Value = 7.8

$ws.Range("C40").Value = 66.1157024793389
$ws.Range("D40").Value = 67.4418604651163
$ws.Range("E40").Value = 66.8

$ws.Range("C41").Value = 242
$ws.Range("D41").Value = 258

# --- Crosstabulasi gender vs awareness ---
$ws.Range("C46").Value = 47.9166666666667
$ws.Range("D46").Value = 49.390243902439
$ws.Range("E46").Value = 48.4

$ws.Range("C47").Value = 52.0833333333333
$ws.Range("D47").Value = 50.609756097561
$ws.Range("E47").Value = 51.6

$ws.Range("C48").Value = 336
$ws.Range("D48").Value = 164

# --- Stasiun TV (single column percentage) ---
$ws.Range("B52").Value = 90.3420523138833
$ws.Range("B53").Value = 72.0321931589537
$ws.Range("B54").Value = 49.8993963782696
$ws.Range("B55").Value = 29.1750503018109
$ws.Range("B56").Value = 1.40845070422535
$ws.Range("B57").Value = 497

# --- Crosstabulasi stasiun TV vs gender ---
$ws.Range("B62").Value = 90.495867768595
$ws.Range("C62").Value = 90.1960784313726
$ws.Range("D62").Value = 90.3420523138833

$ws.Range("B63").Value = 71.0743801652893
$ws.Range("C63").Value = 72.9411764705882
$ws.Range("D63").Value = 72.0321931589537

$ws.Range("B64").Value = 51.2396694214876
$ws.Range("C64").Value = 48.6274509803922
$ws.Range("D64").Value = 49.8993963782696

$ws.Range("B65").Value = 29.7520661157025
$ws.Range("C65").Value = 28.6274509803922
$ws.Range("D65").Value = 29.1750503018109

$ws.Range("B66").Value = 2.06611570247934
$ws.Range("C66").Value = 0.784313725490196
$ws.Range("D66").Value = 1.40845070422535

$ws.Range("B67").Value = 242
$ws.Range("C67").Value = 255
$ws.Range("D67").Value = 497

# --- Crosstabulasi stasiun TV vs gender (A/B variant, usia order) ---
$ws.Range("C73").Value = 66.1
$ws.Range("D73").Value = 67.4
$ws.Range("E73").Value = 66.8

$ws.Range("C74").Value = 11.2
$ws.Range("D74").Value = 15.1
$ws.Range("E74").Value = 13.2

# Rows 75/76 swap age-group labels ("16 - 20 th" and "26 - 30 th")
$ws.Range("B75").Value = "16 - 20 th"
$ws.Range("C75").Value = 9.9
$ws.Range("D75").Value = 7.4
$ws.Range("E75").Value = 8.6

$ws.Range("B76").Value = "26 - 30 th"
$ws.Range("C76").Value = 9.9
$ws.Range("D76").Value = 5.8
$ws.Range("E76").Value = 7.8

$ws.Range("C77").Value = 2.9
$ws.Range("D77").Value = 4.3
$ws.Range("E77").Value = 3.6

$ws.Range("C78").Value = 242
$ws.Range("D78").Value = 258
